$d = $word.ActiveDocument

$replacements = @(
    @("84×96=8064", "70×55=3850"),
    @("43×57=2451", "91×36=3276"),
    @("97×21=2037", "92×15=1380"),
    @("13×28=364", "15×17=255"),
    @("60×92=5520", "54×92=4968"),
    @("86×49=4214", "21×49=1029"),
    @("77×90=6930", "87×68=5916"),
    @("61×47=2867", "79×20=1580"),
    @("84×81=6804", "31×19=589"),
    @("31×53=1643", "60×94=5640"),
    @("30×88=2640", "69×53=3657"),
    @("13×59=767", "90×25=2250"),
    @("92×48=4416", "56×84=4704"),
    @("53×89=4717", "99×19=1881"),
    @("32×13=416", "11×89=979"),
    @("35×26=910", "26×53=1378"),
    @("82×89=7298", "92×72=6624"),
    @("21×29=609", "22×21=462"),
    @("43×53=2279", "17×59=1003"),
    @("35×80=2800", "22×43=946"),
    @("65×81=5265", "47×15=705"),
    @("25×73=1825", "56×78=4368"),
    @("73×65=4745", "48×61=2928"),
    @("28×27=756", "91×72=6552"),
    @("64×13=832", "27×21=567")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
